# "Add files via upload" — the underlying OOXML diff for this commit is a
# pure PowerPoint open/save round-trip (re-numbered sldMasterId/sldId/
# sldLayoutId values, dropped cached <a:fld> display text, added
# <a:endParaRPr> run-less paragraph marks, stripped p14:creationId extLst
# blocks, font PANOSE/theme-family bookkeeping, default-value
# normalization such as <a:stretch/> -> <a:stretch><a:fillRect/></a:stretch>,
# and relationship-id renumbering). None of the actual slide text, shapes,
# images, slide order or slide count changed between the two states.
#
# Those internal bookkeeping fields are not exposed anywhere on the
# PowerPoint object model (there is no property for a slide's numeric
# SlideID, a layout's numeric id, a font's PANOSE string, or the theme
# family GUID), so they cannot be (and should not be) poked at through
# COM automation. Touching the date/slide-number placeholders through
# TextRange would actually do active harm: it would bake the cached
# "02.03.2021" / "<#>" text in as a literal run and destroy the live
# <a:fld> field, which is not what happened upstream.
#
# The faithful COM-level action that matches this commit is therefore
# simply to touch the open presentation without altering any of its
# visible content.

$p = $ppt.ActivePresentation

# Re-save in place; exercises the COM surface without perturbing any
# slide text, shape, picture or ordering.
$p.Save() | Out-Null
